# Commit: "Began reorganizing packages and fixing the GUI"
#
# Adds a new "NotRequired" input column (J) to Sheet1 with a handful of
# sample values, and renames the lone "OK" value in column E (row 12) to
# "lonely".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header for the "NotRequired" field.
$ws.Range("J1").Value = "NotRequired"

# Sample data for the new column (mix of reused strings + new numbers/text).
$ws.Range("J2").Value = "Hello"
$ws.Range("J3").Value = 10
$ws.Range("J4").Value = "World"
$ws.Range("J5").Value = 10.222
$ws.Range("J6").Value = "EnumVal1"

# Rename the stray "OK" value to "lonely".
$ws.Range("E12").Value = "lonely"

# Move the active selection down to E13 (matches the refreshed view state).
$ws.Range("E13").Select() | Out-Null
